$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Caracore"
$ws.Range("B5").Value = "Desenvolvedor Python"
$ws.Range("C5").Value = "suporte@caracore.com.br"
$ws.Range("D5").Value = "2025-07-17 09:32:01"
$ws.Range("F5").Value = "Enviado - Teste"
$ws.Range("G5").Value = "Teste de envio manual"
$ws.Range("H5").Value = 0
